# open file delegate was changed, icon for open file dialog was added
# This script reproduces the edits to the settings table:
#  - adds "substrate_shape" / "substrate_radius" parameter rows after "magnetron_y"
#  - adds "omega_s_max" / "omega_p_max" parameter rows after "NR_max"
#  - updates R_min / R_max / NR_min / x0_1 values
#  - appends a new "rotation_type" parameter row at the end of the table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: insert two new rows at row 8 (before "Длина подложки") for the new
# "substrate_shape" and "substrate_radius" parameters.
# ---------------------------------------------------------------------------
$ws.Range("A8:A9").EntireRow.Insert()

# fix up column-A number formatting / border (it should look like the other
# index cells, e.g. row 4) since freshly inserted rows lose the border style
$ws.Cells.Item(4,1).Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(8,1).Value2 = 6
$ws.Cells.Item(8,2).Value2 = "Форма подложки"
$ws.Cells.Item(8,3).Value2 = "substrate_shape"
$ws.Cells.Item(8,4).Value2 = "'Circle"
$ws.Cells.Item(8,6).Value2 = "cases=['Circle', 'Rectangle']"
$ws.Cells.Item(8,7).Value2 = "model"

$ws.Cells.Item(9,1).Value2 = 7
$ws.Cells.Item(9,2).Value2 = "Радиус подложки"
$ws.Cells.Item(9,3).Value2 = "substrate_radius"
$ws.Cells.Item(9,4).Value2 = 50
# row 9 needs the same quote-prefix cell style as row 7 (magnetron_y) / row 6
$ws.Cells.Item(7,4).Copy()
$ws.Cells.Item(9,4).PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Cells.Item(9,4).Value2 = 50
$ws.Cells.Item(9,5).Value2 = "мм"
$ws.Cells.Item(9,6).Value2 = "+float"
$ws.Cells.Item(9,7).Value2 = "model"

# ---------------------------------------------------------------------------
# Step 2: update R_min / R_max / NR_min values (rows shifted down by 2 after
# the insertion above: R_min row 24->26, R_max row 25->27, NR_min row 28->30)
# ---------------------------------------------------------------------------
$ws.Cells.Item(26,4).Value2 = 70
$ws.Cells.Item(27,4).Value2 = 95
$ws.Cells.Item(30,4).Value2 = 0.01

# ---------------------------------------------------------------------------
# Step 3: insert two new rows before the "Начальное приближение R" (x0_1) row
# (currently row 32) for the new "omega_s_max" / "omega_p_max" parameters.
# ---------------------------------------------------------------------------
$ws.Range("A32:A33").EntireRow.Insert()

$ws.Cells.Item(29,1).Copy()
$ws.Range("A32:A33").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(32,1).Value2 = 30
$ws.Cells.Item(32,2).Value2 = "Макс. угловая скорость солнца"
$ws.Cells.Item(32,3).Value2 = "omega_s_max"
$ws.Cells.Item(32,4).Value2 = 100
$ws.Cells.Item(32,5).Value2 = "оборотов/мин"
$ws.Cells.Item(32,6).Value2 = "+float"
$ws.Cells.Item(32,7).Value2 = "model"

$ws.Cells.Item(33,1).Value2 = 31
$ws.Cells.Item(33,2).Value2 = "Макс. угловая скорость планеты"
$ws.Cells.Item(33,3).Value2 = "omega_p_max"
$ws.Cells.Item(33,4).Value2 = 100
$ws.Cells.Item(33,5).Value2 = "оборотов/мин"
$ws.Cells.Item(33,6).Value2 = "+float"
$ws.Cells.Item(33,7).Value2 = "model"

# ---------------------------------------------------------------------------
# Step 4: update x0_1 ("Начальное приближение R") value, now at row 34
# ---------------------------------------------------------------------------
$ws.Cells.Item(34,4).Value2 = 85

# ---------------------------------------------------------------------------
# Step 5: append a new "rotation_type" row at the end of the table (row 46)
# ---------------------------------------------------------------------------
$ws.Cells.Item(46,1).Value2 = 44
$ws.Cells.Item(45,1).Copy()
$ws.Cells.Item(46,1).PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Cells.Item(46,1).Value2 = 44
$ws.Cells.Item(46,2).Value2 = "Тип вращения"
$ws.Cells.Item(46,3).Value2 = "rotation_type"
$ws.Cells.Item(46,4).Value2 = "Planet"
$ws.Cells.Item(46,6).Value2 = "cases=['Planet', 'Solar']"
$ws.Cells.Item(46,7).Value2 = "model"
$ws.Cells.Item(46,8).Value2 = "Тип вращения подложкодержателя. "

# ---------------------------------------------------------------------------
# Step 6: cosmetic - column E width + view state
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 13.6328125

$ws.Range("E31").Select()
$excel.ActiveWindow.ScrollRow = 19
